$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UEPS")

# Income Statement section
$ws.Range("D20").Value = 50400    # Total Other Income/Expenses Net
$ws.Range("D21").Value = 144800   # Earnings Before Interest And Taxes
$ws.Range("D23").Value = 100400   # Income Before Tax
$ws.Range("D24").Value = 48300    # Income Tax Expense
$ws.Range("D26").Value = 52000    # Income After Tax
$ws.Range("D27").Value = 63600    # Net Income From Continuing Ops
$ws.Range("D32").Value = -50400   # Other Items
$ws.Range("D33").Value = 63300    # Net Income
$ws.Range("D35").Value = 63300    # Net Income Applicable To Common Shares

# Balance Sheet section
$ws.Range("D62").Value = 192900   # Other Liabilities
$ws.Range("E62").Value = 121600
$ws.Range("F62").Value = 122600

$ws.Range("D66").Value = 576800   # Total Liabilities
$ws.Range("E66").Value = 853200
$ws.Range("F66").Value = 770500

$ws.Range("D76").Value = 642500   # Total Stockholder Equity
$ws.Range("E76").Value = 597600
$ws.Range("F76").Value = 493000

# Cash Flow Statement section
$ws.Range("D81").Value = 63300    # Net Income
$ws.Range("D89").Value = 132300   # Total Cash Flow From Operating Activities
$ws.Range("D94").Value = 180700   # Total Cash Flows From Investing Activities
